# BalanceUs - Resumo: refresh the "Resumo dos Dados" summary with corrected
# figures for the Unimed convenio, split Unimed into "Unimed" / "Unimed PF",
# and recompute the trailing "Geral" (grand total) block, which now moves
# down four rows to make room for the new Unimed PF sub-block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force every value to be written as literal text (shared string),
# matching the source workbook where every data cell -- even the numeric-
# looking ones -- is typed as a string. A leading apostrophe is Excel's
# "treat as text" marker; ClearFormats() afterwards drops the transient
# quote-prefix cell style that the apostrophe trick leaves behind so cell
# styling is left untouched.
function Set-Text {
    param($addr, $val)
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).ClearFormats()
}

# --- Corrected figures for the existing "Unimed" rows (19-22) -------------
Set-Text "E19" "219"
Set-Text "F19" "17278,13"
Set-Text "G19" "1411,06"

Set-Text "E21" "305"
Set-Text "F21" "25514,24"
Set-Text "G21" "2340,26"

Set-Text "E22" "20"
Set-Text "F22" "1570,61"

# --- Rows 23-26 become the new "Unimed PF" sub-block -----------------------
Set-Text "A23" "2020-06-26"
Set-Text "B23" "2020-07-24"
Set-Text "C23" "Unimed PF"
Set-Text "D23" "Gerusa"
Set-Text "E23" "68"
Set-Text "F23" "5994,24"
Set-Text "G23" "458,68"
Set-Text "H23" "0,00"
Set-Text "I23" "0,00"

Set-Text "A24" "2020-06-26"
Set-Text "B24" "2020-07-24"
Set-Text "C24" "Unimed PF"
Set-Text "D24" "Laurise"
Set-Text "E24" "0"
Set-Text "F24" "0,00"
Set-Text "G24" "0,00"
Set-Text "H24" "0,00"
Set-Text "I24" "0,00"

Set-Text "A25" "2020-06-26"
Set-Text "B25" "2020-07-24"
Set-Text "C25" "Unimed PF"
Set-Text "D25" "Valéria"
Set-Text "E25" "58"
Set-Text "F25" "5083,59"
Set-Text "G25" "407,75"
Set-Text "H25" "0,00"
Set-Text "I25" "0,00"

Set-Text "A26" "2020-06-26"
Set-Text "B26" "2020-07-24"
Set-Text "C26" "Unimed PF"
Set-Text "D26" "Procedimentos"
Set-Text "E26" "2"
Set-Text "F26" "128,38"

# --- New rows 27-30: the "Geral" (grand total) block, now recomputed and
#     shifted down to make room for the Unimed PF rows above --------------
$ws.Rows("27:30").Insert()

Set-Text "A27" ""
Set-Text "B27" ""
Set-Text "C27" "Geral"
Set-Text "D27" "Gerusa"
Set-Text "E27" "294"
Set-Text "F27" "35502,29"

Set-Text "A28" ""
Set-Text "B28" ""
Set-Text "C28" "Geral"
Set-Text "D28" "Laurise"
Set-Text "E28" "0"
Set-Text "F28" "0,00"

Set-Text "A29" ""
Set-Text "B29" ""
Set-Text "C29" "Geral"
Set-Text "D29" "Valéria"
Set-Text "E29" "377"
Set-Text "F29" "53003,55"

Set-Text "A30" ""
Set-Text "B30" ""
Set-Text "C30" "Geral"
Set-Text "D30" "Procedimentos"
Set-Text "E30" "22"
Set-Text "F30" "3141,22"
